$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.259.68"
$ws.Range("E2").Value = "  +0.84%  "

$ws.Range("D3").Value = "1.920.06"
$ws.Range("E3").Value = "  +0.44%  "

$ws.Range("E4").Value = "  +0.35%  "

$ws.Range("D5").Value = "'0.8070"
$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("D6").Value = "'244.44"
$ws.Range("E6").Value = "  +0.96%  "

$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("D8").Value = "'0.3254"
$ws.Range("E8").Value = "  +2.71%  "

$ws.Range("D9").Value = "'26.95"
$ws.Range("E9").Value = "  +1.66%  "

$ws.Range("D10").Value = "'0.07268"
$ws.Range("E10").Value = "  +4.95%  "

$ws.Range("D11").Value = "'0.7889"
$ws.Range("E11").Value = "  +6.41%  "

$ws.Range("D12").Value = "'0.08087"
$ws.Range("E12").Value = "  +1.11%  "

$ws.Range("D13").Value = "1.929.87"
$ws.Range("E13").Value = "  +1.17%  "

$ws.Range("D14").Value = "'5.408"
$ws.Range("E14").Value = "  +3.95%  "

$ws.Range("D15").Value = "'93.95"
$ws.Range("E15").Value = "  +0.93%  "

$ws.Range("D16").Value = "30.264.27"
$ws.Range("E16").Value = "  +0.92%  "

$ws.Range("E17").Value = "  +1.51%  "

$ws.Range("D18").Value = "'6.071"
$ws.Range("E18").Value = "  +3.07%  "

$ws.Range("D19").Value = "'250.11"
$ws.Range("E19").Value = "  +1.59%  "

$ws.Range("D20").Value = "'0.000007856"
$ws.Range("E20").Value = "  +1.34%  "

$ws.Range("D21").Value = "2.178.71"
$ws.Range("E21").Value = "  +0.73%  "

$ws.Range("D22").Value = "'8.215"
$ws.Range("E22").Value = "  +20.01%  "

$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("D25").Value = "'0.1652"
$ws.Range("E25").Value = "  +16.04%  "

$ws.Range("D26").Value = "'9.494"
$ws.Range("E26").Value = "  +2.94%  "

$ws.Range("D27").Value = "'167.94"
$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("E28").Value = "  +0.40%  "

$ws.Range("D29").Value = "'2.159"
$ws.Range("E29").Value = "  +5.78%  "

$ws.Range("D30").Value = "'1.394"
$ws.Range("E30").Value = "  +2.18%  "

$ws.Range("D31").Value = "'1.552"
$ws.Range("E31").Value = "  +2.29%  "

$ws.Range("D32").Value = "'4.366"
$ws.Range("E32").Value = "  +1.15%  "

$ws.Range("D33").Value = "'0.05737"
$ws.Range("E33").Value = "  +4.44%  "

$ws.Range("D34").Value = "'4.146"
$ws.Range("E34").Value = "  +1.34%  "

$ws.Range("D35").Value = "'1.298"
$ws.Range("E35").Value = "  +2.16%  "

$ws.Range("D36").Value = "'0.7492"
$ws.Range("E36").Value = "  +2.40%  "

$ws.Range("D37").Value = "'1.006"
$ws.Range("E37").Value = "  +0.75%  "

$ws.Range("D38").Value = "'2.731"
$ws.Range("E38").Value = "  +0.46%  "

$ws.Range("D39").Value = "'0.01961"
$ws.Range("E39").Value = "  +1.89%  "

$ws.Range("D40").Value = "'2.821"
$ws.Range("E40").Value = "  +1.33%  "

$ws.Range("D41").Value = "'0.4543"
$ws.Range("E41").Value = "  +2.55%  "

$ws.Range("D42").Value = "'74.26"
$ws.Range("E42").Value = "  +2.26%  "

$ws.Range("D43").Value = "'6.000"
$ws.Range("E43").Value = "  -2.63%  "

$ws.Range("D44").Value = "'0.8543"
$ws.Range("E44").Value = "  +2.10%  "

$ws.Range("D45").Value = "'1.931"
$ws.Range("E45").Value = "  +2.74%  "

$ws.Range("E46").Value = "  +0.11%  "

$ws.Range("D47").Value = "'103.52"
$ws.Range("E47").Value = "  +3.03%  "

$ws.Range("D48").Value = "1.033.50"
$ws.Range("E48").Value = "  +4.91%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'10.08"
$ws.Range("E49").Value = "  +3.45%  "

$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "'7.625"
$ws.Range("E50").Value = "  +0.71%  "

$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "'3.091"
$ws.Range("E51").Value = "  +11.11%  "
